# "fix party name and substrings"
#
# 1. The Respondent party placeholder was wrong ({{case.client}} should be
#    {{case.respondent}}).
# 2. The slice syntax on the SSN/driver's-license substrings was backwards
#    ([-1:3] should be [-3:] to get "last three characters").
# 3. A comma was added after "postjudgment interest" in the attorney-fees
#    paragraph.
$d = $word.ActiveDocument

# --- 1. Respondent, {{case.client}}, files this Original Answer... ---
$p2 = $d.Paragraphs.Item(2).Range
$p2.Find.Execute("{{case.client}}", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "{{case.respondent}}", 2)

# --- 2 & 3. SSN / driver's-license substring slices ---
$p3 = $d.Paragraphs.Item(3).Range
$p3.Find.Execute("ssn[-1:3]", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "ssn[-3:]", 2)

$p3 = $d.Paragraphs.Item(3).Range
$p3.Find.Execute("dl_number[-1:3]", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "dl_number[-3:]", 2)

# --- 4. "...postjudgment interest should be..." -> "...postjudgment interest, should be..." ---
# A comma is inserted right after "postjudgment interest". Word's hidden
# _GoBack bookmark (which tracks the location of the most recent edit) moves
# along with that edit, so it is relocated from its old spot (just before
# the second "{{case.me.name}}") to just after the new comma.
$text = $d.Content.Text

$favorIdx = $text.IndexOf("in favor of {{case.")
$oldGoBackPos = $favorIdx + ("in favor of {{case.").Length
$d.Bookmarks.Add("ZZ_keep_boundary", $d.Range($oldGoBackPos, $oldGoBackPos))

$interestIdx = $text.IndexOf("postjudgment interest")
$commaPos = $interestIdx + ("postjudgment interest").Length
$d.Range($commaPos, $commaPos).InsertBefore(",")

$d.Bookmarks.Add("_GoBack", $d.Range($commaPos + 1, $commaPos + 1))
$d.Bookmarks.Item("ZZ_keep_boundary").Delete()
